$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

$ws.Range("C2").Value = 4230.41705472264
$ws.Range("F2").Value = -0.265976874331424

$ws.Range("C3").Value = 4395.31256292677
$ws.Range("F3").Value = 6.60466869391803

$ws.Range("C4").Value = 4406.60344171618
$ws.Range("F4").Value = 117.075121508101

$ws.Range("C5").Value = 4370.40690777762
$ws.Range("F5").Value = 115.56693366971

$ws.Range("C6").Value = 3700.16386776256
$ws.Range("F6").Value = 87.6401403356652

$ws.Range("B7").Value = 975.205470885634
$ws.Range("C7").Value = 1576.02970777409
$ws.Range("F7").Value = 19.4493550786857

$ws.Range("C8").Value = 1518.64846132174
$ws.Range("F8").Value = 17.7032493821693

$ws.Range("C9").Value = 4945.68451789699
$ws.Range("F9").Value = 108.679362717993

$ws.Range("C10").Value = 5263.70679100814
$ws.Range("F10").Value = 121.932959767669

$ws.Range("C11").Value = 5248.25683702003
$ws.Range("F11").Value = 121.289211684831

$ws.Range("C12").Value = 5216.86343966236
$ws.Range("F12").Value = 120.731574331105

$ws.Range("C13").Value = 4496.11918008676
$ws.Range("F13").Value = 90.7005635154549

$ws.Range("C14").Value = 2239.45247202791
$ws.Range("F14").Value = 12.6799069154369

$ws.Range("C15").Value = 2241.49468832503
$ws.Range("F15").Value = 12.2253378092816
